# Insert a new weekly price record for "Vega Modelo de Temuco" (Apio)
# as row 153, pushing the existing rows 153-184 down to 154-185.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 153..184 down by one to make room for the new record.
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the new observation.
$ws.Cells.Item(153, 1).Value  = 10
$ws.Cells.Item(153, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(153, 3).Value  = 'La Araucanía'
$ws.Cells.Item(153, 4).Value  = 44476
$ws.Cells.Item(153, 5).Value  = 9
$ws.Cells.Item(153, 6).Value  = 100112017
$ws.Cells.Item(153, 7).Value  = 'Apio'
$ws.Cells.Item(153, 8).Value  = 'Americana (o)'
$ws.Cells.Item(153, 9).Value  = 'Primera'
$ws.Cells.Item(153, 10).Value = 100
$ws.Cells.Item(153, 11).Value = 9000
$ws.Cells.Item(153, 12).Value = 9000
$ws.Cells.Item(153, 13).Value = 9000
$ws.Cells.Item(153, 14).Value = '$/docena de matas'
$ws.Cells.Item(153, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(153, 16).Value = 1500
$ws.Cells.Item(153, 17).Value = 6
$ws.Cells.Item(153, 18).Value = 'Hortaliza'

# Keep the date column formatted like the rest of column D.
$ws.Cells.Item(153, 4).NumberFormat = $ws.Cells.Item(154, 4).NumberFormat
